# Implemented skills distribution analysis
# Append 8 new job-application rows (47-54) mirroring the existing data layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A = 46; B = 43990; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Not found";   Skills = "Python HTML CSS Javascript Java GIT REST " },
    @{ A = 47; B = 43990; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Not found";   Skills = "Python HTML CSS Javascript Java C++ GIT REST " },
    @{ A = 48; B = 43990; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Toronto, ON"; Skills = "Python HTML CSS Javascript Java GIT REST " },
    @{ A = 49; B = 43990; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Not found";   Skills = "Python HTML CSS Javascript Java GIT REST " },
    @{ A = 50; B = 43991; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Not found";   Skills = "Python HTML CSS Javascript Java GIT REST " },
    @{ A = 51; B = 43991; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Toronto, ON"; Skills = "Python HTML CSS Javascript Java C++ GIT REST " },
    @{ A = 52; B = 43991; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Toronto, ON"; Skills = "Python HTML CSS Javascript Java GIT REST " },
    @{ A = 53; B = 43991; Title = "Web Developer Intern"; Company = "Pathcore"; Location = "Toronto, ON"; Skills = "Python HTML CSS Javascript Java GIT REST " }
)

$url = "https://ca.indeed.com/viewjob?jk=45e7adfb4d34664e&tk=1e8k9749r0gc1000&from=serp&vjs=3"

$startRow = 47
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data.A

    $ws.Cells.Item($row, 2).Value = $data.B

    $ws.Cells.Item($row, 3).Value = $data.Title
    $ws.Cells.Item($row, 4).Value = $data.Company
    $ws.Cells.Item($row, 5).Value = $data.Location
    $ws.Cells.Item($row, 6).Value = $data.Skills

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $url, "", "", $url)
}
